$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: NIK (A), Plat (B), Nama (C), Password (D)
# A3 is a long numeric-looking value that must stay text (like A2), so force
# text entry the same way the sheet's author typed it (leading apostrophe),
# then drop the quote-prefix formatting bit so no style index lingers.
$ws.Range("A3").Value = "'1234456278949533"
$ws.Range("A3").ClearFormats()

# B3 is blank, just like B2 - an explicit empty text value.
$ws.Range("B3").Value = "'"
$ws.Range("B3").ClearFormats()

$ws.Range("C3").Value = "Rahma"
$ws.Range("D3").Value = "Akun_rahma21"
